# The sheet originally held 13 columns (A:M) of Katalon-style element
# locator data. The commit ("Generated by Katalon AI") trims the sheet
# down to just the "input_Name" column (which was column I), dropping
# every other column (A-H and J-M).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the trailing columns (J:M) first, then the leading ones (A:H),
# so the surviving column I slides down into column A without needing
# any extra re-indexing math.
$ws.Columns("J:M").Delete()
$ws.Columns("A:H").Delete()
